$d = $word.ActiveDocument

# Remove the unnecessary "Blah blalh blahbbbh" text (along with the
# spell-check proofErr markers surrounding "blalh") while leaving the
# paragraph's bookmarkStart/bookmarkEnd (_GoBack) intact.
$d.Content.Find.Execute("Blah blalh blahbbbh", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
